$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G50").Value = 71.09999999999999
$ws.Range("E51").Value = 91.40000000000001
$ws.Range("E52").Value = 89.7
$ws.Range("B54").Value = 80.90000000000001
$ws.Range("D54").Value = 89.5
$ws.Range("F54").Value = 82
$ws.Range("E57").Value = 86.5
$ws.Range("G60").Value = 76.90000000000001
$ws.Range("G62").Value = 83.7
$ws.Range("H68").Value = 97.3
$ws.Range("E71").Value = 101.2
$ws.Range("G71").Value = 100.5
$ws.Range("G72").Value = 100.8
$ws.Range("E74").Value = 98.8
$ws.Range("G76").Value = 102.6
$ws.Range("B78").Value = 103
$ws.Range("F78").Value = 100.9
$ws.Range("G78").Value = 104.3
$ws.Range("G79").Value = 104.6
$ws.Range("C80").Value = 100.7
$ws.Range("G80").Value = 103.9
$ws.Range("J80").Value = 104.8
$ws.Range("G81").Value = 103.7
$ws.Range("C82").Value = 102.8
$ws.Range("E82").Value = 98.8
$ws.Range("F82").Value = 106.1
$ws.Range("D83").Value = 100.3
$ws.Range("G83").Value = 104.6
$ws.Range("I83").Value = 105.5
$ws.Range("J83").Value = 106
$ws.Range("B84").Value = 105.7
$ws.Range("C84").Value = 101.4
$ws.Range("E84").Value = 97.2
$ws.Range("F84").Value = 108.5
$ws.Range("G84").Value = 104.6
$ws.Range("H84").Value = 108.9
$ws.Range("I84").Value = 105.6
$ws.Range("J84").Value = 106.7
$ws.Range("B85").Value = 106
$ws.Range("G85").Value = 106.2
$ws.Range("D86").Value = 89.3
$ws.Range("E86").Value = 97.5
$ws.Range("G86").Value = 107.7
$ws.Range("D87").Value = 96.90000000000001
$ws.Range("F87").Value = 105.3
$ws.Range("G87").Value = 110.3
$ws.Range("I87").Value = 106.1
$ws.Range("J87").Value = 107.3
$ws.Range("B88").Value = 108.1
$ws.Range("C88").Value = 103.7
$ws.Range("F88").Value = 107.9
$ws.Range("G88").Value = 109
$ws.Range("I88").Value = 107.9
$ws.Range("J88").Value = 108.8
$ws.Range("G89").Value = 110.6
$ws.Range("D90").Value = 105.9
$ws.Range("E90").Value = 101.3
$ws.Range("F90").Value = 110.3
$ws.Range("G90").Value = 112.7
$ws.Range("B91").Value = 111.4
$ws.Range("C91").Value = 106.2
$ws.Range("D91").Value = 101.7
$ws.Range("E91").Value = 105
$ws.Range("F91").Value = 110.5
$ws.Range("G91").Value = 114.6
$ws.Range("J91").Value = 112.4
$ws.Range("B92").Value = 110.9
$ws.Range("C92").Value = 104.8
$ws.Range("D92").Value = 100.6
$ws.Range("E92").Value = 103.5
$ws.Range("F92").Value = 108.9
$ws.Range("G92").Value = 113.4
$ws.Range("H92").Value = 114.1
$ws.Range("I92").Value = 110.5
$ws.Range("J92").Value = 112
$ws.Range("E93").Value = 104.2
$ws.Range("G93").Value = 113.3
$ws.Range("B94").Value = 112
$ws.Range("D94").Value = 96.59999999999999
$ws.Range("F94").Value = 111.2
$ws.Range("G94").Value = 114.4
$ws.Range("I94").Value = 111.7
$ws.Range("B95").Value = 113.3
$ws.Range("C95").Value = 105.5
$ws.Range("D95").Value = 97.8
$ws.Range("E95").Value = 104.5
$ws.Range("F95").Value = 112.1
$ws.Range("G95").Value = 117
$ws.Range("H95").Value = 117.3
$ws.Range("I95").Value = 112.9
$ws.Range("J95").Value = 114.9
$ws.Range("B96").Value = 113.9
$ws.Range("C96").Value = 107.4
$ws.Range("D96").Value = 102.4
$ws.Range("E96").Value = 105.2
$ws.Range("F96").Value = 112.9
$ws.Range("G96").Value = 114.5
$ws.Range("H96").Value = 117.9
$ws.Range("I96").Value = 113.6
$ws.Range("J96").Value = 115.1
$ws.Range("B97").Value = 109.6
$ws.Range("E97").Value = 101.8
$ws.Range("F97").Value = 111.1
$ws.Range("G97").Value = 107.8
$ws.Range("J97").Value = 110.4
$ws.Range("C98").Value = 106.3
$ws.Range("D98").Value = 101.8
$ws.Range("F98").Value = 112.4
$ws.Range("G98").Value = 111.3
$ws.Range("H98").Value = 115.8
$ws.Range("J98").Value = 112.8
$ws.Range("B99").Value = 97.8
$ws.Range("C99").Value = 97.2
$ws.Range("D99").Value = 100.8
$ws.Range("E99").Value = 92.5
$ws.Range("F99").Value = 98.09999999999999
$ws.Range("G99").Value = 94.7
$ws.Range("H99").Value = 99.40000000000001
$ws.Range("J99").Value = 97.3
$ws.Range("B100").Value = 102.8
$ws.Range("C100").Value = 98
$ws.Range("D100").Value = 101.9
$ws.Range("E100").Value = 98.7
$ws.Range("F100").Value = 94.40000000000001
$ws.Range("G100").Value = 113.1
$ws.Range("H100").Value = 103.1
$ws.Range("I100").Value = 102.2
$ws.Range("J100").Value = 102.8
$ws.Range("B101").Value = 109.5
$ws.Range("C101").Value = 103.8
$ws.Range("F101").Value = 106
$ws.Range("G101").Value = 122.6
$ws.Range("J101").Value = 110.6
$ws.Range("B102").Value = 113.3
$ws.Range("C102").Value = 104.7
$ws.Range("D102").Value = 100.8
$ws.Range("E102").Value = 106.3
$ws.Range("F102").Value = 107.1
$ws.Range("G102").Value = 125.3
$ws.Range("H102").Value = 116.4
$ws.Range("I102").Value = 112.7
$ws.Range("J102").Value = 114.9
$ws.Range("B103").Value = 114.6
$ws.Range("C103").Value = 106.7
$ws.Range("D103").Value = 103.8
$ws.Range("E103").Value = 107.3
$ws.Range("F103").Value = 108.9
$ws.Range("G103").Value = 136.3
$ws.Range("H103").Value = 114.1
$ws.Range("I103").Value = 113.4
$ws.Range("J103").Value = 116
